$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label fix ---
$ws.Range("A1").Value = "Nutrient"

# --- Encoding / wording fixes on nutrient display names ---
$ws.Range("A5").Value  = "Calcium (mg)"
$ws.Range("A6").Value  = "Folic acid (µg DFE)"
$ws.Range("A7").Value  = "Iron absorbed (mg)"
$ws.Range("A11").Value = "Pantothenic acid (mg)"
$ws.Range("A12").Value = "Vitamin A (μg RAE)"
$ws.Range("A16").Value = "Vitamin B12 (µg)"

# --- Row 12 (Vitamin A): mark as having an Upper_limits linkage and add the
#     food-composition nutrient name used for the upper-limit match ---
$ws.Range("G12").Value = "RetinolActivityEquivalent"
$ws.Range("D12").Value = $true

# --- Header row: center-align ---
$ws.Range("A1:H1").HorizontalAlignment = -4108

# --- Selection cursor left where the author finished working ---
$ws.Range("J12").Select()
